$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.962.42"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "'3.524.32"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'586.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "'178.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'3.523.12"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").Value = "'6.93"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "'0.426"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").Value = "'4.137.67"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'30.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.19%  "
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "'66.950.34"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "'3.524.98"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "'6.09"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").Value = "'14.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "'383.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'72.52"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'24.76"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").Value = "'5.94"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.42%  "
$ws.Range("D33").Value = "'2.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").Value = "'1.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").Value = "'7.31"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").Value = "'1.59"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").Value = "'30.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +14.78%  "
$ws.Range("D39").Value = "'161.30"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "'0.899"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").Value = "'6.65"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").Value = "'4.55"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  -8.08%  "
$ws.Range("D45").Value = "'2.739.51"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").Value = "'0.0708"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "'25.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.33%  "
$ws.Range("D48").Value = "'40.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "'0.0299"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'324.45"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("D51").Value = "'1.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.98%  "
